# "ooutput update 2025 august"
#
# Refresh the generated IG-publisher output for the
# mindfulness-audit-format StructureDefinition:
#   - canonical URL moved from the old GitHub-shorthand location to the
#     2rdoc.pt IG site
#   - regeneration date bumped to the August 2025 run
#   - the audit-formats ValueSet binding now also lives under 2rdoc.pt
#   - the "Elements" sheet's auto-fit column widths shrink slightly to
#     match the refreshed content

$wb = $excel.ActiveWorkbook

$wsMeta     = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------

# URL (row 2)
$wsMeta.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/mindfulness-audit-format"

# Date (row 8)
$wsMeta.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# --- Elements sheet --------------------------------------------------------

# Binding Value Set for the audit-formats required binding (Z6)
$wsElements.Range("Z6").Value = "https://2rdoc.pt/fhir/ValueSet/audit-formats"

# The publisher re-runs AutoFit over the whole sheet on every generation,
# so the column widths drift slightly whenever any cell content changes.
# Reproduce the refreshed best-fit widths column by column.
$colWidths = [ordered]@{
    1  = 15.666666666666666
    2  = 15.666666666666666
    3  = 9.0
    4  = 6.166666666666667
    5  = 4.5
    6  = 3.1666666666666665
    7  = 3.5
    8  = 11.833333333333334
    9  = 9.666666666666666
    10 = 19.833333333333332
    11 = 7.5
    12 = 99.83333333333333
    13 = 99.83333333333333
    14 = 99.83333333333333
    15 = 11.5
    16 = 19.833333333333332
    17 = 19.833333333333332
    18 = 19.833333333333332
    19 = 19.833333333333332
    20 = 7.0
    21 = 12.833333333333334
    22 = 13.166666666666666
    23 = 14.166666666666666
    24 = 13.833333333333334
    25 = 16.166666666666668
    26 = 35.166666666666664
    27 = 4.166666666666667
    28 = 17.166666666666668
    29 = 33.666666666666664
    30 = 12.666666666666666
    31 = 10.5
    32 = 14.166666666666666
    33 = 7.333333333333333
    34 = 7.666666666666667
    35 = 99.83333333333333
    37 = 18.666666666666668
}

foreach ($col in $colWidths.Keys) {
    $wsElements.Columns.Item($col).ColumnWidth = $colWidths[$col]
}

# These columns (C, D, AE, AF, AG) were hidden before the edit and stay
# hidden afterwards - re-assert it since re-writing ColumnWidth above
# touches the <col> element.
$hiddenCols = @(3, 4, 31, 32, 33)
foreach ($col in $hiddenCols) {
    $wsElements.Columns.Item($col).Hidden = $true
}
